$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Project start date moved from 08-Mar-2017 to 03-Aug-2017.
#    All other dates in the plan are driven off this via formulas
#    (B6=C4, C6=B6+6, B7=B6+7, C7=B7+6, ... shared down to row 20),
#    so they recalc automatically.
$ws.Range("C4").Value = 42950

# 2) Rewrite the text of the first three weekly task cells (content was
#    reshuffled between weeks 1-3 and the testing plan was added).
$ws.Range("D6").Value = "Phân tích yêu cầu và lập usecase tổng quát - Team`nThảo luận, thống nhất usecase, kiến trúc, chia task - Team`nTạo project repository trên Github - Nguyễn Thanh Nga"

$ws.Range("D7").Value = "Tạo sơ đồ thực thể, sơ đồ thực thể quan hệ, sơ đồ lớp, view flow - Team`nPhân tích usecase, vẽ lược đồ trình tự (sequence diagram) và lược đồ hoạt động (activity diagram) - Theo usecase được phân công`nLập kế hoạch kiểm thử - Nguyễn Thanh Nga"

$ws.Range("D8").Value = "Tạo khung project, tạo database - Nguyễn Thanh Nga`nThảo luận, hoàn thiện các lược đồ - Team`nVẽ mockup cho các view,  lập bảng Test-case - Theo usecase được phân công"

# D8 now wraps to three lines, so the row needs to grow to match D7.
$ws.Rows.Item(8).RowHeight = 62.4

# 3) Fill in the remaining weeks (previously-empty plan cells) with the
#    coding / testing plan, matching the existing D6:D8 look (bold 12pt
#    Times New Roman). Some entries wrap to two lines, others fit on one.
$ws.Range("D9").Value = "Coding - Theo usecase được phân công"
$ws.Range("D10").Value = "Coding - Theo usecase được phân công"
$ws.Range("D11").Value = "Coding - Theo usecase được phân công"
$ws.Range("D12").Value = "Coding - Theo usecase được phân công"
$ws.Range("D13").Value = "Tiến hành kiểm thử chức năng, hiệu chỉnh - Theo usecase được phân công"
$ws.Range("D14").Value = "Tiến hành kiểm thử chức năng, hiệu chỉnh - Theo usecase được phân công"
$ws.Range("D15").Value = "Kiểm thử tích hợp , hiệu chỉnh - Team"
$ws.Range("D16").Value = "Kiểm thử tích hợp , hiệu chỉnh - Team"
$ws.Range("D17").Value = "Kiểm thử hệ thống, hiệu chỉnh - Team"
$ws.Range("D18").Value = "Kiểm thử hệ thống, hiệu chỉnh,  lập tài liệu hướng dẫn sử dụng - Team"
$ws.Range("D19").Value = "Báo cáo đề tài, triển khai hệ thống - Team"
$ws.Range("D20").Value = "Báo cáo đề tài, triển khai hệ thống - Team"

$newPlanRange = $ws.Range("D9:D20")
$newPlanRange.Font.Bold = $true
$newPlanRange.Font.Size = 12
$newPlanRange.Font.Name = "Times New Roman"

$ws.Range("D9:D12").WrapText = $false
$ws.Range("D13:D14").WrapText = $true
$ws.Range("D15:D16").WrapText = $false
$ws.Range("D17:D18").WrapText = $true
$ws.Range("D19:D20").WrapText = $false

# 4) Touch the (new) row below the table - formatting spills one row
#    further down once the plan is complete.
$ws.Range("D21").Font.Size = 12
$ws.Range("D21").Font.Name = "Times New Roman"
$ws.Range("D21").VerticalAlignment = -4108

# 5) Move the on-screen selection to reflect where the user ended up
#    working (week 3 area).
$ws.Range("D12").Select()
